# Adapt column header formatting to respective input file names:
#   "<name>_old" -> "<name>_FV2410"
#   "<name>_new" -> "<name>_FV2504"
# and expose the header row as a proper Excel Table (ListObject), plus
# freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- 1. Rename the header cells (row 1) -----------------------------------
$oldHeaders = @(
    "Segmentname_old",
    "Segmentgruppe_old",
    "Segment_old",
    "Datenelement_old",
    "Segment ID_old",
    "Code_old",
    "Qualifier_old",
    "Beschreibung_old",
    "Bedingungsausdruck_old",
    "Bedingung_old"
)

$newHeaders = @(
    "Segmentname_new",
    "Segmentgruppe_new",
    "Segment_new",
    "Datenelement_new",
    "Segment ID_new",
    "Code_new",
    "Qualifier_new",
    "Beschreibung_new",
    "Bedingungsausdruck_new",
    "Bedingung_new"
)

$fv2410Headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

$fv2504Headers = @(
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

# Columns A..J hold the "_old" -> "_FV2410" headers
for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2410Headers[$i]
}

# Column K holds "diff" and is left untouched.

# Columns L..U hold the "_new" -> "_FV2504" headers
for ($i = 0; $i -lt $newHeaders.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2504Headers[$i]
}

# --- 2. Turn the data range into an Excel Table (ListObject) --------------
$dataRange = $ws.Range("A1:U59")
$table = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $dataRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$table.Name = "Table1"
$table.TableStyle = "TableStyleMedium9"

# --- 3. Freeze the header row ----------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

$wb.Save()
